# Update Sheets per scheduled runner refresh (Brynhildr_Profits)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 6901.0835
$ws.Range("I32").Value = 2699
$ws.Range("J32").Value = 9902.571
$ws.Range("K32").Value = 2699
$ws.Range("L32").Value = 9902.571
$ws.Range("M32").Value = -2373
$ws.Range("N32").Value = -10554.571
$ws.Range("H98").Value = 1343.0667
$ws.Range("I98").Value = 1153.2858
$ws.Range("K98").Value = 1153.2858
$ws.Range("M98").Value = 344.7141999999999
$ws.Range("H111").Value = 168.85715
$ws.Range("J111").Value = 122.5
$ws.Range("L111").Value = 367.5
$ws.Range("N111").Value = -6501.5
$ws.Range("H113").Value = 4238.1113
$ws.Range("I113").Value = 3624.75
$ws.Range("J113").Value = 4728.8
$ws.Range("K113").Value = 3624.75
$ws.Range("L113").Value = 4728.8
$ws.Range("M113").Value = -370.75
$ws.Range("N113").Value = -11236.8
$ws.Range("H122").Value = 1343.0667
$ws.Range("I122").Value = 1153.2858
$ws.Range("K122").Value = 3459.8574
$ws.Range("M122").Value = -1009.8574
$ws.Range("H125").Value = 3084
$ws.Range("I125").Value = 3084
$ws.Range("J125").Value = 0
$ws.Range("K125").Value = 27756
$ws.Range("L125").Value = 0
$ws.Range("M125").Value = -25296
$ws.Range("N125").ClearContents()
$ws.Range("H127").Value = 5729.5
$ws.Range("I127").Value = 3594.25
$ws.Range("K127").Value = 10782.75
$ws.Range("M127").Value = -5822.75
$ws.Range("H137").Value = 23815598
$ws.Range("J137").Value = 26227.5
$ws.Range("L137").Value = 78682.5
$ws.Range("N137").Value = -83782.5
$ws.Range("H138").Value = 3919.3372
$ws.Range("I138").Value = 6390.1816
$ws.Range("K138").Value = 19170.5448
$ws.Range("M138").Value = -14030.5448
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 1714.4736
$ws.Range("I110").Value = 1698.6111
$ws.Range("K110").Value = 1698.6111
$ws.Range("M110").Value = 346.3888999999999
$ws.Range("H122").Value = 1901.5
$ws.Range("I122").Value = 1781.1818
$ws.Range("K122").Value = 5343.5454
$ws.Range("M122").Value = -2893.5454
$ws.Range("H133").Value = 66504.664
$ws.Range("J133").Value = 66504.664
$ws.Range("L133").Value = 66504.664
$ws.Range("N133").Value = -71564.664
$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").ClearContents()
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 504.5
$ws.Range("I22").Value = 422.66666
$ws.Range("J22").Value = 750
$ws.Range("K22").Value = 422.66666
$ws.Range("L22").Value = 750
$ws.Range("M22").Value = -249.66666
$ws.Range("N22").Value = -1096
$ws.Range("H132").Value = 76704.5
$ws.Range("J132").Value = 76704.5
$ws.Range("L132").Value = 76704.5
$ws.Range("N132").Value = -86824.5
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 47226.52
$ws.Range("I16").Value = 11394
$ws.Range("J16").Value = 80302.69500000001
$ws.Range("K16").Value = 11394
$ws.Range("L16").Value = 80302.69500000001
$ws.Range("M16").Value = -11107
$ws.Range("N16").Value = -80876.69500000001
$ws.Range("H31").Value = 1881262
$ws.Range("I31").Value = 3290976.8
$ws.Range("J31").Value = 1642.25
$ws.Range("K31").Value = 3290976.8
$ws.Range("L31").Value = 1642.25
$ws.Range("M31").Value = -3290681.8
$ws.Range("N31").Value = -2232.25
$ws.Range("H34").Value = 1881262
$ws.Range("I34").Value = 3290976.8
$ws.Range("J34").Value = 1642.25
$ws.Range("K34").Value = 3290976.8
$ws.Range("L34").Value = 1642.25
$ws.Range("M34").Value = -3290774.8
$ws.Range("N34").Value = -2046.25
$ws.Range("H86").Value = 47297.863
$ws.Range("I86").Value = 81355.414
$ws.Range("J86").Value = 6428.8
$ws.Range("K86").Value = 81355.414
$ws.Range("L86").Value = 6428.8
$ws.Range("M86").Value = -80232.414
$ws.Range("N86").Value = -8674.799999999999
$ws.Range("H89").Value = 47297.863
$ws.Range("I89").Value = 81355.414
$ws.Range("J89").Value = 6428.8
$ws.Range("K89").Value = 406777.07
$ws.Range("L89").Value = 32144
$ws.Range("M89").Value = -401161.07
$ws.Range("N89").Value = -43376
$ws.Range("H98").Value = 66979.5
$ws.Range("J98").Value = 66979.5
$ws.Range("L98").Value = 66979.5
$ws.Range("N98").Value = -71471.5
$ws.Range("H99").Value = 11760.305
$ws.Range("I99").Value = 17281.643
$ws.Range("J99").Value = 3171.5557
$ws.Range("K99").Value = 17281.643
$ws.Range("L99").Value = 3171.5557
$ws.Range("M99").Value = -15783.643
$ws.Range("N99").Value = -6167.5557
$ws.Range("H107").Value = 1204.9166
$ws.Range("I107").Value = 807.625
$ws.Range("K107").Value = 807.625
$ws.Range("M107").Value = 1112.375
$ws.Range("H113").Value = 47226.52
$ws.Range("I113").Value = 11394
$ws.Range("J113").Value = 80302.69500000001
$ws.Range("K113").Value = 11394
$ws.Range("L113").Value = 80302.69500000001
$ws.Range("M113").Value = -9224
$ws.Range("N113").Value = -84642.69500000001
$ws.Range("H117").Value = 85000
$ws.Range("J117").Value = 85000
$ws.Range("L117").Value = 85000
$ws.Range("N117").Value = -94178
$ws.Range("H122").Value = 8173.8
$ws.Range("I122").Value = 1930.3158
$ws.Range("J122").Value = 126800
$ws.Range("K122").Value = 5790.9474
$ws.Range("L122").Value = 380400
$ws.Range("M122").Value = -3340.9474
$ws.Range("N122").Value = -385300
$ws.Range("H126").Value = 11760.305
$ws.Range("I126").Value = 17281.643
$ws.Range("J126").Value = 3171.5557
$ws.Range("K126").Value = 51844.929
$ws.Range("L126").Value = 9514.667099999999
$ws.Range("M126").Value = -49374.929
$ws.Range("N126").Value = -14454.6671
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1735318.8
$ws.Range("I5").Value = 1701438.8
$ws.Range("J5").Value = 1790047.9
$ws.Range("K5").Value = 5104316.4
$ws.Range("L5").Value = 5370143.699999999
$ws.Range("M5").Value = -5104204.4
$ws.Range("N5").Value = -5370367.699999999
$ws.Range("H8").Value = 238.375
$ws.Range("I8").Value = 238.375
$ws.Range("K8").Value = 715.125
$ws.Range("M8").Value = -576.125
$ws.Range("H121").Value = 18342.5
$ws.Range("J121").Value = 22475.924
$ws.Range("L121").Value = 67427.772
$ws.Range("N121").Value = -70047.772
$ws.Range("H135").Value = 1735318.8
$ws.Range("I135").Value = 1701438.8
$ws.Range("J135").Value = 1790047.9
$ws.Range("K135").Value = 15312949.2
$ws.Range("L135").Value = 16110431.1
$ws.Range("M135").Value = -15310414.2
$ws.Range("N135").Value = -16115501.1
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H31").Value = 2682.8572
$ws.Range("I31").Value = 2296.6667
$ws.Range("K31").Value = 2296.6667
$ws.Range("M31").Value = -2004.6667
$ws.Range("H37").Value = 2682.8572
$ws.Range("I37").Value = 2296.6667
$ws.Range("K37").Value = 2296.6667
$ws.Range("M37").Value = -2019.6667
$ws.Range("H70").Value = 18681.346
$ws.Range("I70").Value = 15777.107
$ws.Range("J70").Value = 100000
$ws.Range("K70").Value = 15777.107
$ws.Range("L70").Value = 100000
$ws.Range("M70").Value = -15507.107
$ws.Range("N70").Value = -100540
$ws.Range("H73").Value = 18681.346
$ws.Range("I73").Value = 15777.107
$ws.Range("J73").Value = 100000
$ws.Range("K73").Value = 15777.107
$ws.Range("L73").Value = 100000
$ws.Range("M73").Value = -14841.107
$ws.Range("N73").Value = -101872
$ws.Range("H122").Value = 3286.125
$ws.Range("I122").Value = 3540.4
$ws.Range("J122").Value = 2862.3333
$ws.Range("K122").Value = 10621.2
$ws.Range("L122").Value = 8586.999899999999
$ws.Range("M122").Value = -8171.200000000001
$ws.Range("N122").Value = -13486.9999
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 3508.2273
$ws.Range("I46").Value = 1188.7778
$ws.Range("K46").Value = 1188.7778
$ws.Range("M46").Value = -1000.7778
$ws.Range("H68").Value = 7677.095
$ws.Range("I68").Value = 8346.786
$ws.Range("J68").Value = 6337.7144
$ws.Range("K68").Value = 8346.786
$ws.Range("L68").Value = 6337.7144
$ws.Range("M68").Value = -7597.786
$ws.Range("N68").Value = -7835.7144
$ws.Range("H71").Value = 7677.095
$ws.Range("I71").Value = 8346.786
$ws.Range("J71").Value = 6337.7144
$ws.Range("K71").Value = 41733.93
$ws.Range("L71").Value = 31688.572
$ws.Range("M71").Value = -37989.93
$ws.Range("N71").Value = -39176.572
$ws.Range("H82").Value = 2919.7273
$ws.Range("I82").Value = 2953.75
$ws.Range("K82").Value = 2953.75
$ws.Range("M82").Value = -2592.75
$ws.Range("H85").Value = 2919.7273
$ws.Range("I85").Value = 2953.75
$ws.Range("K85").Value = 2953.75
$ws.Range("M85").Value = -1705.75
$ws.Range("H114").Value = 40000
$ws.Range("J114").Value = 40000
$ws.Range("L114").Value = 40000
$ws.Range("N114").Value = -48678
$ws.Range("H122").Value = 7700.8
$ws.Range("I122").Value = 4504
$ws.Range("J122").Value = 8500
$ws.Range("K122").Value = 13512
$ws.Range("L122").Value = 25500
$ws.Range("M122").Value = -11062
$ws.Range("N122").Value = -30400
$ws.Range("H134").Value = 90000
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 90000
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 90000
$ws.Range("M134").ClearContents()
$ws.Range("N134").Value = -100140
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H16").Value = 68760
$ws.Range("J16").Value = 68760
$ws.Range("L16").Value = 68760
$ws.Range("N16").Value = -69344
$ws.Range("H38").Value = 2500
$ws.Range("J38").Value = 0
$ws.Range("L38").Value = 0
$ws.Range("N38").ClearContents()
$ws.Range("H48").Value = 0
$ws.Range("J48").Value = 0
$ws.Range("L48").Value = 0
$ws.Range("N48").ClearContents()
$ws.Range("H54").Value = 0
$ws.Range("J54").Value = 0
$ws.Range("L54").Value = 0
$ws.Range("N54").ClearContents()
